# Update NATMI ligand-receptor edge metrics with newly computed TPM-based values.
# (Ligand/Receptor average & total expression, derived specificities, and edge weights
# for rows 2-9 of the Zp3-Egfr sheet.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.08261866666666666
$ws.Range("H2").Value = 0.247856
$ws.Range("I2").Value = 0.5806915650061265
$ws.Range("J2").Value = 0.5806915650061265
$ws.Range("M2").Value = 1.370876333333333
$ws.Range("N2").Value = 4.112629
$ws.Range("O2").Value = 0.01103063309339269
$ws.Range("P2").Value = 0.01103063309339269
$ws.Range("Q2").Value = 0.1132599748248889
$ws.Range("R2").Value = 1.019339773424
$ws.Range("S2").Value = 0.006405395594010574
$ws.Range("T2").Value = 0.006405395594010573
$ws.Range("G3").Value = 0.08261866666666666
$ws.Range("H3").Value = 0.247856
$ws.Range("I3").Value = 0.5806915650061265
$ws.Range("J3").Value = 0.5806915650061265
$ws.Range("O3").Value = 0.7476219244149905
$ws.Range("P3").Value = 0.7476219244149904
$ws.Range("Q3").Value = 7.676408019454222
$ws.Range("R3").Value = 69.087672175088
$ws.Range("S3").Value = 0.4341377453214329
$ws.Range("T3").Value = 0.4341377453214328
$ws.Range("G4").Value = 0.08261866666666666
$ws.Range("H4").Value = 0.247856
$ws.Range("I4").Value = 0.5806915650061265
$ws.Range("J4").Value = 0.5806915650061265
$ws.Range("M4").Value = 29.718484
$ws.Range("N4").Value = 89.155452
$ws.Range("O4").Value = 0.2391271080585153
$ws.Range("P4").Value = 0.2391271080585153
$ws.Range("Q4").Value = 2.455301523434667
$ws.Range("R4").Value = 22.097713710912
$ws.Range("S4").Value = 0.1388590946138884
$ws.Range("T4").Value = 0.1388590946138884
$ws.Range("G5").Value = 0.08261866666666666
$ws.Range("H5").Value = 0.247856
$ws.Range("I5").Value = 0.5806915650061265
$ws.Range("J5").Value = 0.5806915650061265
$ws.Range("M5").Value = 0.275941
$ws.Range("N5").Value = 0.827823
$ws.Range("O5").Value = 0.002220334433101459
$ws.Range("P5").Value = 0.002220334433101458
$ws.Range("Q5").Value = 0.02279787749866666
$ws.Range("R5").Value = 0.205180897488
$ws.Range("S5").Value = 0.001289329476794677
$ws.Range("T5").Value = 0.001289329476794677
$ws.Range("I6").Value = 0.4193084349938734
$ws.Range("J6").Value = 0.4193084349938734
$ws.Range("M6").Value = 1.370876333333333
$ws.Range("N6").Value = 4.112629
$ws.Range("O6").Value = 0.01103063309339269
$ws.Range("P6").Value = 0.01103063309339269
$ws.Range("Q6").Value = 0.08178328333522222
$ws.Range("R6").Value = 0.736049550017
$ws.Range("S6").Value = 0.004625237499382119
$ws.Range("T6").Value = 0.004625237499382118
$ws.Range("I7").Value = 0.4193084349938734
$ws.Range("J7").Value = 0.4193084349938734
$ws.Range("O7").Value = 0.7476219244149905
$ws.Range("P7").Value = 0.7476219244149904
$ws.Range("S7").Value = 0.3134841790935576
$ws.Range("T7").Value = 0.3134841790935575
$ws.Range("I8").Value = 0.4193084349938734
$ws.Range("J8").Value = 0.4193084349938734
$ws.Range("M8").Value = 29.718484
$ws.Range("N8").Value = 89.155452
$ws.Range("O8").Value = 0.2391271080585153
$ws.Range("P8").Value = 0.2391271080585153
$ws.Range("Q8").Value = 1.772935412310667
$ws.Range("R8").Value = 15.956418710796
$ws.Range("S8").Value = 0.1002680134446269
$ws.Range("T8").Value = 0.1002680134446269
$ws.Range("I9").Value = 0.4193084349938734
$ws.Range("J9").Value = 0.4193084349938734
$ws.Range("M9").Value = 0.275941
$ws.Range("N9").Value = 0.827823
$ws.Range("O9").Value = 0.002220334433101459
$ws.Range("P9").Value = 0.002220334433101458
$ws.Range("Q9").Value = 0.01646199619766666
$ws.Range("R9").Value = 0.148157965779
$ws.Range("S9").Value = 0.0009310049563067817
$ws.Range("T9").Value = 0.0009310049563067815
